# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new value for column F
$exhibitUpdates = @{
    7  = 2390
    9  = 222
    11 = 2508
    16 = 132
    18 = 9352
    19 = 58
    20 = 7264
    21 = 11828
    27 = 2647
    29 = 204
    30 = 2600
    31 = 794
    32 = 50
    34 = 977
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型": row -> new value for column F
$allTypesUpdates = @{
    11 = 2390
    14 = 222
    15 = 2508
    21 = 132
    23 = 9352
    25 = 7264
    26 = 11828
    34 = 2647
    38 = 204
    39 = 50
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}
